# Update crypto price/volume data (coinranking.com scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.918.80"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.606.13"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'210.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.485"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.86%  "
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("D10").Value = "'17.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("D11").Value = "'0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "1.828.22"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("D13").Value = "1.601.76"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("E15").Value = "  -3.28%  "
$ws.Range("D16").Value = "25.886.34"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "'61.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'190.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("D22").Value = "'9.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").Value = "'5.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").Value = "'0.129"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'142.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("E28").Value = "  -3.34%  "
$ws.Range("D29").Value = "'14.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").Value = "'0.0470"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("E33").Value = "  -4.20%  "
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "1.117.27"
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("E38").Value = "  -6.94%  "
$ws.Range("D39").Value = "'0.0151"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("D40").Value = "'0.501"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.96%  "
$ws.Range("D41").Value = "'96.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("D42").Value = "1.740.72"
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").Value = "'0.746"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.83%  "
$ws.Range("D44").Value = "'5.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.23%  "
$ws.Range("D45").Value = "0.0₆0114"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "'53.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.38%  "
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("E48").Value = "  -2.77%  "
$ws.Range("D49").Value = "'0.411"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "'7.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.37%  "
